# Rename "在庫" (stock table) sheet to "在庫履歴" (stock history),
# and rework its schema: the id formula reference switches from 工場 to
# 原価, and four new PK/field rows are inserted (品目ID, 時刻, 在庫単位,
# 在庫数量) before the existing 在庫金額 row.

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("在庫")
$ws.Name = "在庫履歴"

# Title / physical table name
$ws.Range("A1").Value = "在庫履歴"
$ws.Range("B1").Value = "stock_history"

# Make room for the new field rows (品目ID, 時刻, 在庫単位, 在庫数量)
# between the existing 工場ID row (4) and 在庫金額 row (old row 5).
$ws.Rows("5:8").Insert()

# Row 4 (工場ID) now pulls its type/length from 原価 instead of 工場
$ws.Range("C4").Formula = "=原価!C3"
$ws.Range("D4").Formula = "=原価!D3"

# Row 5: 品目ID (new PK field)
$ws.Range("A5").Value = "品目ID"
$ws.Range("B5").WrapText = $true
$ws.Range("B5").Value = "s_i_id"
$ws.Range("C5").Formula = "=原価!C4"
$ws.Range("D5").Formula = "=原価!D4"
$ws.Range("F5").Value = "PK"

# Row 6: 時刻 (new PK field, literal data type/length)
$ws.Range("A6").Value = "時刻"
$ws.Range("B6").WrapText = $true
$ws.Range("B6").Value = "s_time"
$ws.Range("C6").Value = "time"
$ws.Range("D6").Value = 6
$ws.Range("F6").Value = "PK"
$ws.Range("G6").Value = "時分秒"

# Row 7: 在庫単位 (new non-key field)
$ws.Range("A7").Value = "在庫単位"
$ws.Range("B7").ClearFormats()
$ws.Range("B7").Value = "s_stock_unit"
$ws.Range("C7").Formula = "=原価!C5"
$ws.Range("D7").Formula = "=原価!D5"
$ws.Range("G7").Value = "原価マスター．在庫単位"

# Row 8: 在庫数量 (new non-key field)
$ws.Range("A8").Value = "在庫数量"
$ws.Range("B8").ClearFormats()
$ws.Range("B8").Value = "s_stock_quantity"
$ws.Range("C8").Formula = "=原価!C6"
$ws.Range("D8").Formula = "=原価!D6"
$ws.Range("E8").Formula = "=原価!E6"
$ws.Range("G8").Value = "原価マスター．在庫数量"

# Row 9 (formerly row 5): 在庫金額 label simplified, formulas/content unchanged
$ws.Range("A9").Value = "在庫金額"

# Restore cursor positions to match the author's final selection
$ws5 = $wb.Worksheets.Item("原価")
$ws5.Activate()
$ws5.Range("D10").Select()

$ws.Activate()
$ws.Range("E12").Select()
